$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.295.43"
$ws.Range("E2").Value = "  +0.79%  "

$ws.Range("D3").Value = "1.863.07"
$ws.Range("E3").Value = "  +0.87%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'0.7024"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("D6").Value = "'237.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "'0.08218"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.84%  "

$ws.Range("D9").Value = "'0.3052"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.17%  "

$ws.Range("D10").Value = "'23.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.21%  "

$ws.Range("D11").Value = "'0.08185"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("D12").Value = "1.866.32"
$ws.Range("E12").Value = "  +1.32%  "

$ws.Range("D13").Value = "'0.7196"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.92%  "

$ws.Range("D14").Value = "'5.194"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.71%  "

$ws.Range("D15").Value = "'89.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("D16").Value = "29.304.73"
$ws.Range("E16").Value = "  +1.10%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "'5.804"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.28%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000007891"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.84%  "

$ws.Range("E19").Value = "  +2.90%  "

$ws.Range("D20").Value = "'237.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.90%  "

$ws.Range("D21").Value = "'1.0000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").Value = "2.107.34"
$ws.Range("E22").Value = "  +1.85%  "

$ws.Range("D23").Value = "'1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.25%  "

$ws.Range("D24").Value = "'7.490"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.85%  "

$ws.Range("D25").Value = "'162.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.92%  "

$ws.Range("D26").Value = "'9.013"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.60%  "

$ws.Range("D27").Value = "'0.1452"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.88%  "

$ws.Range("E28").Value = "  +0.85%  "

$ws.Range("E29").Value = "  +2.84%  "

$ws.Range("E30").Value = "  +3.72%  "

$ws.Range("D31").Value = "'4.441"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.30%  "

$ws.Range("E32").Value = "  -0.64%  "

$ws.Range("D33").Value = "'4.073"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.67%  "

$ws.Range("D34").Value = "'0.05223"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.19%  "

$ws.Range("D35").Value = "'1.175"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.89%  "

$ws.Range("D36").Value = "'0.7075"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("D37").Value = "'1.003"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.13%  "

$ws.Range("D38").Value = "'2.661"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.60%  "

$ws.Range("E39").Value = "  -0.54%  "

$ws.Range("D40").Value = "'2.721"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.76%  "

$ws.Range("D41").Value = "1.151.03"
$ws.Range("E41").Value = "  +8.47%  "

$ws.Range("D42").Value = "'0.9206"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.34%  "

$ws.Range("D43").Value = "'5.972"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("D44").Value = "'0.4291"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").Value = "'71.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.41%  "

$ws.Range("D46").Value = "'0.9994"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").Value = "'103.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.24%  "

$ws.Range("D48").Value = "'1.781"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.61%  "

$ws.Range("D49").Value = "2.005.59"
$ws.Range("E49").Value = "  +1.15%  "

$ws.Range("D50").Value = "'9.210"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.49%  "

$ws.Range("D51").Value = "'6.996"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.88%  "
